# #5: property boat&car done
# Rebuild the "汽車" (car) sheet: it gained a "capacity" column and the
# full set of normalized metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that the
# other property sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- header row (row 1) --------------------------------------------------
# Existing header cells B1/D1/E1/F1 already hold the right text (owner,
# register_date, register_reason) - only B1/C1 change meaning, and new
# trailing header cells H1:N1 need to be created. All header cells use the
# bold + bordered + centered style already used in row 1.
$headerCells = @(
    @("B1", "name"),
    @("C1", "capacity"),
    @("D1", "owner"),
    @("E1", "register_date"),
    @("F1", "register_reason"),
    @("G1", "acquire_value"),
    @("H1", "property_category"),
    @("I1", "category"),
    @("J1", "date"),
    @("K1", "legislator_name"),
    @("L1", "legislator_id"),
    @("M1", "source_file"),
    @("N1", "index")
)
foreach ($pair in $headerCells) {
    $cell = $ws.Range($pair[0])
    $cell.Value = $pair[1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# ---- data rows (2-4) ------------------------------------------------------
# B column used to hold the car name directly ("小客車") in both the
# header-looking row 1 and each data row; now it's purely data and the
# extra metadata columns (H:N) are appended to every row.
$rows = @(
    @{ row = 2; idx = 51; capacity = 1781; owner = "郭綺雯"; register_date = "96年12月28日"; acquire_value = 800000 },
    @{ row = 3; idx = 52; capacity = 4799; owner = "郭綺雯"; register_date = "100年01月11曰"; acquire_value = 1450000 },
    @{ row = 4; idx = 53; capacity = 2461; owner = "紀國棟"; register_date = "99年09月14日"; acquire_value = 1200000 }
)

foreach ($r in $rows) {
    $row = $r.row
    $ws.Range("B$row").Value = "小客車"
    $ws.Range("C$row").Value = $r.capacity
    $ws.Range("D$row").Value = $r.owner
    $ws.Range("E$row").Value = $r.register_date
    $ws.Range("F$row").Value = "買賣"
    $ws.Range("G$row").Value = $r.acquire_value
    $ws.Range("H$row").Value = "land"
    $ws.Range("I$row").Value = "normal"

    # "date" (J) holds the literal text "2011-12-20" in every other sheet
    # (shared-string, not a real date). Force text entry so Excel's COM
    # layer doesn't silently coerce it into a date serial number, then
    # drop the number-format override so the cell is left with plain
    # default formatting (matching every other un-styled data cell).
    $dcell = $ws.Range("J$row")
    $dcell.NumberFormat = "@"
    $dcell.Value = "2011-12-20"
    $dcell.ClearFormats()

    $ws.Range("K$row").Value = "紀國棟"
    $ws.Range("L$row").Value = 918
    $ws.Range("M$row").Value = "tmp5e8b1"
    $ws.Range("N$row").Value = $r.idx
}
